# "moving model training to rubi"
# Converts the per-sample free-text colony/streak counts (col B, some rows
# also had stray D/E concentration columns) into a consistent layout:
#   col A = sample id, col B = count description (text) or raw count (number),
#   col C = computed concentration (count * dilution, mostly as a formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 had stray per-prep concentration numbers in D2/E2 - replace with the
# same style formula used elsewhere in column C, then drop the old D/E data
# that is no longer part of the table (columns D:E are unused everywhere else).
$ws.Range("D1:E81").ClearContents()

$ws.Range("C2").Formula  = "=14*100*1000"
$ws.Range("C3").Formula  = "=22*100*1000"
$ws.Range("C4").Formula  = "=4*100*1000"
$ws.Range("C5").Formula  = "=12*100*1000"
$ws.Range("C6").Formula  = "=4*100*1000"
$ws.Range("C7").Formula  = "=12*100*1000"
$ws.Range("C8").Formula  = "=15*100*1000"
$ws.Range("C9").Formula  = "=15*100*1000"

# Row 10: record the raw colony count as a number instead of the
# "325 colonies" text description (C10 already holds the literal count).
$ws.Range("B10").Value = 325

# Row 11: now has a count description + computed concentration.
$ws.Range("B11").Value   = "9 confluent streaks"
$ws.Range("C11").Formula = "=9*100*1000"

# Row 12: raw numeric colony count, like row 10.
$ws.Range("B12").Value = 425

# Row 13
$ws.Range("B13").Value   = "5 confluent streaks"
$ws.Range("C13").Formula = "=5*100*1000"

# Row 14
$ws.Range("B14").Value   = "9 confluent streaks"
$ws.Range("C14").Formula = "=9*100*1000"

# Row 15: raw numeric count with a relative (B15-based) formula, matching
# the pattern used by other numeric-count rows (27-29, 36-47, etc).
$ws.Range("B15").Value   = 105
$ws.Range("C15").Formula = "=B15*1000"

# Row 16
$ws.Range("B16").Value   = "11 confluent streaks"
$ws.Range("C16").Formula = "=11*100*1000"

# Row 17
$ws.Range("B17").Value   = "7 confluent streaks"
$ws.Range("C17").Formula = "=7*100*1000"

# Row 18
$ws.Range("B18").Value   = 100
$ws.Range("C18").Formula = "=B18*1000"

# Row 19
$ws.Range("B19").Value   = "11 confluent streaks"
$ws.Range("C19").Formula = "=11*100*1000"

# Row 20
$ws.Range("B20").Value   = "17 confluent streaks"
$ws.Range("C20").Formula = "=17*100*1000"

# Row 21
$ws.Range("B21").Value   = "7 confluent streaks"
$ws.Range("C21").Formula = "=7*100*1000"

# Row 22
$ws.Range("B22").Value   = "5 confluent streaks"
$ws.Range("C22").Formula = "=5*100*1000"

# Row 23
$ws.Range("B23").Value   = "22 confluent streaks"
$ws.Range("C23").Formula = "=22*100*1000"

# Row 24
$ws.Range("B24").Value   = "5 confluent streaks"
$ws.Range("C24").Formula = "=5*100*1000"

# Row 25
$ws.Range("B25").Value   = "3 confluent streaks"
$ws.Range("C25").Formula = "=3*100*1000"

# Row 26
$ws.Range("B26").Value   = "4 confluent streaks"
$ws.Range("C26").Formula = "=4*100*1000"

# Rows 27-29 already had their B/C values populated; leave as-is.

# Row 30: description already correct, just add the computed concentration.
$ws.Range("C30").Formula = "=4*100*1000"

# Row 31: description text re-confirmed, concentration added.
$ws.Range("B31").Value   = "3 confluent streaks"
$ws.Range("C31").Formula = "=3*100*1000"

# Row 32 unchanged (no count recorded).

# Row 33
$ws.Range("B33").Value   = "21 confluent streaks"
$ws.Range("C33").Formula = "=21*100*1000"

# Row 34
$ws.Range("B34").Value   = "7 confluent streaks"
$ws.Range("C34").Formula = "=7*100*1000"

# Row 35: description already correct, just add the computed concentration.
$ws.Range("C35").Formula = "=12*100*1000"

# Restore selection/scroll state to match the refreshed table extent.
$ws.Range("B1:C81").Select()
